$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$tb = $s.Shapes.AddTextbox(1, 116.19, 450.88, 524.47, 29.08)
$tb.Name = "TextBox 1"
$tb.TextFrame.WordWrap = $false
$tr = $tb.TextFrame.TextRange
$tr.Text = "https://docs.python.org/3/library/pdb.html#debugger-commands "
